$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = "Hartmut"

# B3 holds a long numeric-looking card number that must stay text;
# prefix with an apostrophe so Excel keeps it as text instead of a number.
$ws.Range("B3").Value = "'2570314725427075"
$ws.Range("C3").Value = "Mohaupt"

$ws.Range("D5").Value = "KONTOSTAND AM 18.08.2025"

$ws.Range("B6").Value = "19.08."
$ws.Range("C6").Value = "20.08."
$ws.Range("D6").Value = "ABSCHLAG STROM Stadtwerke Rosenheim 79689669"
$ws.Range("E6").Value = "86,95-"

$ws.Range("B7").Value = "21.08."
$ws.Range("C7").Value = "22.08."
$ws.Range("D7").Value = "MCDONALDS Demmin"
$ws.Range("E7").Value = "9,62-"

$ws.Range("B8").Value = "23.08."
$ws.Range("C8").Value = "24.08."
$ws.Range("D8").Value = "RECHNUNG VODAFONE GMBH 16868108"
$ws.Range("E8").Value = "38,11-"

$ws.Range("B9").Value = "26.08."
$ws.Range("C9").Value = "27.08."
$ws.Range("D9").Value = "BURGER KING Mainburg"
$ws.Range("E9").Value = "25,51-"

$ws.Range("D12").Value = "KONTOSTAND AM 29.08.2025"
$ws.Range("E12").Value = "160,19-"

$ws.Range("C13").Value = "IHR NAECHSTER ABRECHNUNGSTERMIN 07.09.2025"
